$d = $word.ActiveDocument

# Locate the very last paragraph in the document body (the empty w:p right
# before the sectPr). We insert a new empty paragraph just ahead of it and
# then pour the new Section II content into that fresh paragraph, so the
# pre-existing trailing empty paragraph is left completely untouched.
$count = $d.Paragraphs.Count
$tailRange = $d.Paragraphs.Item($count).Range
$tailRange.Collapse(1)             # wdCollapseStart
$tailRange.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($count)

$sectionXml = @"
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
      </w:pPr>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>Section II: Mitigation, Continuity, and Disasters</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve">Hi-Tech operates on a finite budget and must prioritize investments into features and services.  On the one hand, the business would like to spend all available resources delivering its core mission, building the best electronic vehicle.  Allocating time and money into other projects might even appear to detract from this </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>mission, and</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> feel like a waste.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:ind w:firstLine="720"/>
      </w:pPr>
      <w:r>
        <w:t>However, investments in other aspects of the organization reduce risk and improve continuity.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">  </w:t>
      </w:r>
      <w:r>
        <w:t>It can be challenging to choose efficiently prioritize risk reduction because it compares an immediate real cost against a hypothetical future expense</w:t>
      </w:r>
      <w:sdt>
        <w:sdtPr>
          <w:id w:val="1310443408"/>
          <w:citation/>
        </w:sdtPr>
        <w:sdtContent>
          <w:r>
            <w:fldChar w:fldCharType="begin"/>
          </w:r>
          <w:r>
            <w:instrText xml:space="preserve"> CITATION Gor15 \l 1033 </w:instrText>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="separate"/>
          </w:r>
          <w:r>
            <w:rPr>
              <w:noProof/>
            </w:rPr>
            <w:t xml:space="preserve"> (Gordon, Loeb, Lucyshyn, &amp; Zhou, 2015)</w:t>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="end"/>
          </w:r>
        </w:sdtContent>
      </w:sdt>
      <w:r>
        <w:t>.  For instance, licensing anti-virus software might cost the business one million dollars a year.  If during that year, the company was lucky and did not encounter any malware, then the insurance was not used.  In contrast, ransomware spreading across the intranet could easily exceed several million dollars</w:t>
      </w:r>
      <w:sdt>
        <w:sdtPr>
          <w:id w:val="1204684966"/>
          <w:citation/>
        </w:sdtPr>
        <w:sdtContent>
          <w:r>
            <w:fldChar w:fldCharType="begin"/>
          </w:r>
          <w:r>
            <w:instrText xml:space="preserve"> CITATION Ast16 \l 1033 </w:instrText>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="separate"/>
          </w:r>
          <w:r>
            <w:rPr>
              <w:noProof/>
            </w:rPr>
            <w:t xml:space="preserve"> (Astani &amp; Ready, 2016)</w:t>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="end"/>
          </w:r>
        </w:sdtContent>
      </w:sdt>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:ind w:firstLine="720"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Along with continuity solutions, the business also requires disaster recovery and response solutions to handle both known and unknown assaults.  </w:t>
      </w:r>
      <w:r>
        <w:t>The threat landscape continues to evolve with adversaries</w:t>
      </w:r>
      <w:r>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> continually gaining leverage through decreasing costs to issue the attack</w:t>
      </w:r>
      <w:r>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> versus </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">the cost to the </w:t>
      </w:r>
      <w:r>
        <w:t>defend</w:t>
      </w:r>
      <w:r>
        <w:t>er</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> protect</w:t>
      </w:r>
      <w:r>
        <w:t>ing</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>their resources</w:t>
      </w:r>
      <w:sdt>
        <w:sdtPr>
          <w:id w:val="962086271"/>
          <w:citation/>
        </w:sdtPr>
        <w:sdtContent>
          <w:r>
            <w:fldChar w:fldCharType="begin"/>
          </w:r>
          <w:r>
            <w:instrText xml:space="preserve"> CITATION Lam16 \l 1033 </w:instrText>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="separate"/>
          </w:r>
          <w:r>
            <w:rPr>
              <w:noProof/>
            </w:rPr>
            <w:t xml:space="preserve"> (Lam, 2016)</w:t>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="end"/>
          </w:r>
        </w:sdtContent>
      </w:sdt>
      <w:r>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">  This asymmetry naturally occurs because the attacker only needs to send packets versus the defender must parse and act on those requests.  Malware authors are also publishing over one million strains every day</w:t>
      </w:r>
      <w:sdt>
        <w:sdtPr>
          <w:id w:val="-461274516"/>
          <w:citation/>
        </w:sdtPr>
        <w:sdtContent>
          <w:r>
            <w:fldChar w:fldCharType="begin"/>
          </w:r>
          <w:r>
            <w:instrText xml:space="preserve"> CITATION Kil171 \l 1033 </w:instrText>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="separate"/>
          </w:r>
          <w:r>
            <w:rPr>
              <w:noProof/>
            </w:rPr>
            <w:t xml:space="preserve"> (Kilgallon, De La Rosa, &amp; Cavazos, 2017)</w:t>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="end"/>
          </w:r>
        </w:sdtContent>
      </w:sdt>
      <w:r>
        <w:t>.  Given the sheer volume, there is a high probability that anti-virus fails to prevent an infection</w:t>
      </w:r>
      <w:r>
        <w:t>.  While these malicious actors can wreak havoc on the corporate network, few forces are as damaging as employee negligence</w:t>
      </w:r>
      <w:sdt>
        <w:sdtPr>
          <w:id w:val="1880124459"/>
          <w:citation/>
        </w:sdtPr>
        <w:sdtContent>
          <w:r>
            <w:fldChar w:fldCharType="begin"/>
          </w:r>
          <w:r>
            <w:instrText xml:space="preserve"> CITATION Val171 \l 1033 </w:instrText>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="separate"/>
          </w:r>
          <w:r>
            <w:rPr>
              <w:noProof/>
            </w:rPr>
            <w:t xml:space="preserve"> (Valiente, 2017)</w:t>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="end"/>
          </w:r>
        </w:sdtContent>
      </w:sdt>
      <w:r>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">  If a support technician issues an erroneous database command, it can cascade into a critical outage.  Even after protecting against these sources, a hurricane or fire can cause irreparable damage.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>

"@

$newPara.Range.InsertXML($sectionXml)

Write-Host "Paragraphs after insert:" $d.Paragraphs.Count
